$wb = $excel.ActiveWorkbook

# --- Sheet1: update "track" policy param2 row (row 4) ---
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Range("B4").Value = 2
$ws1.Range("C4").Value = ""

# --- Sheet2: insert two new columns (track_min, track_max) before increment_size ---
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Columns("E:F").Insert()
$ws2.Range("E1").Value = "track_min"
$ws2.Range("F1").Value = "track_max"
$ws2.Range("E2").Value = 0
$ws2.Range("F2").Value = 4
$ws2.Range("G2").Value = 0.1

# --- Update selections on every sheet to match the saved state ---
$ws1.Range("C4").Select()

$ws4 = $wb.Worksheets.Item("Sheet4")
$ws4.Range("E13").Select()

$ws3 = $wb.Worksheets.Item("Sheet3")
$ws3.Range("C2").Select()

# Sheet2 becomes the active/selected tab last
$ws2.Range("F13").Select()
